$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = '[''have'', ''deliver'', ''work'', ''drive'', ''help'', ''connect'', ''support'', ''know'', ''build'', ''bring'', ''leverage'', ''focus'', ''manage'', ''root'', ''want'', ''own'', ''create'', ''transform'', ''interact'', ''lead'', ''oversee'', ''do'', ''enable'', ''solutionsmarkete'', ''solve'', ''sell'', ''expect'', ''suit'', ''make'', ''allow'', ''understand'', ''matter'', ''respond'', ''unlock'', ''empower'', ''respect'', ''give'', ''personalize'', ''seek'', ''elevate'', ''set'', ''end'', ''automate'', ''decisione'', ''ensure'', ''meet'', ''include'', ''continue'', ''pay'', ''mentor'', ''join'', ''grow'', ''maintain'', ''structure'', ''align'', ''engage'', ''inspire'', ''act'', ''implement'', ''operate'', ''market'', ''serve'', ''follow'', ''operationsadvertise'', ''base'', ''accelerate'', ''resultsqualificationsrequire'', ''bringyou'', ''develop'', ''expand'', ''kpi'', ''provide'', ''learn'', ''improve'', ''take'', ''pull'', ''demonstrate'', ''staff'', ''achieve'', ''position'', ''triage'', ''prioritize'', ''compete'', ''pace'', ''prove'', ''produce'', ''write'', ''authorize'', ''travel'', ''suspend'', ''vary'']'
$ws.Range("I3").Value = '[''work'', ''drive'', ''lead'', ''solve'', ''take'', ''provide'', ''prioritize'', ''communicate'', ''manage'', ''relate'', ''integrate'', ''address'', ''cross'', ''responsibilitiesidentify'', ''benefit'', ''validate'', ''utilize'', ''review'', ''models.contribute'', ''feature'', ''help'', ''scope'', ''timeline'', ''deal'', ''deliver'', ''enable'', ''scale'', ''define'', ''create'', ''compete'', ''translate'', ''technologiesoutstande'', ''do'', ''follow'', ''perform'', ''impact'', ''develop'', ''implement'', ''preferred.applie'', ''boost'', ''inspire'', ''maximize'', ''dig'', ''think'', ''pace'', ''break'', ''defend'', ''debate'', ''use'', ''base'']'
$ws.Range("I4").Value = '[''identify'', ''include'', ''build'', ''work'', ''improve'', ''communicate'', ''seek'', ''use'', ''apply'', ''solve'', ''drive'', ''protect'', ''support'', ''become'', ''develop'', ''increase'', ''provide'', ''grow'', ''analyze'', ''automate'', ''innovate'', ''look'', ''desire'', ''model'', ''discover'', ''cut'', ''design'', ''answer'', ''resolution'', ''define'', ''deliver'', ''unlock'', ''detect'', ''source'', ''conduct'', ''query'', ''offline'', ''casescollaborate'', ''enhance'', ''understand'', ''know'', ''use.maintain'', ''diverse'', ''presentationsabout'', ''leverage'', ''make'', ''effortless'', ''relate'', ''script'', ''write'', ''present'', ''datapreferre'', ''demonstrate'', ''exceed'', ''prioritize'', ''sacrifice'', ''dive'', ''commit'', ''discriminate'', ''like'', ''request'', ''visit'', ''a1985797'']'
$ws.Range("I5").Value = '[''require'', ''include'', ''learn'', ''leverage'', ''develop'', ''validate'', ''communicate'', ''make'', ''receive'', ''need'', ''follow'', ''provide'', ''allow'', ''commit'', ''set'', ''evaluate'', ''understand'', ''regard'', ''establish'', ''procedures.require'', ''relate'', ''toolsexperience'', ''solve'', ''write'', ''work'', ''contribute'', ''focus'', ''improve'', ''associate'', ''have'', ''lack'', ''protect'', ''informationpreferre'', ''manage'', ''perform'', ''build'', ''covid'', ''undergo'', ''wear'', ''cover'', ''do'', ''approve'', ''come'', ''put'', ''progress'', ''offer'', ''test'', ''support'', ''use'', ''interview'', ''call'', ''enhance'', ''connect'', ''gain'', ''pertain'', ''select'', ''check'', ''ensure'', ''miss'', ''invite'', ''participate'', ''listen'', ''anticipate'', ''take'', ''review'', ''inform'', ''move'', ''interviews.schedule'']'
$ws.Range("I6").Value = '[''build'', ''solve'', ''apply'', ''learn'', ''’re'', ''advance'', ''obsess'', ''make'', ''change'', ''consume'', ''partner'', ''lead'', ''transform'', ''use'', ''excite'', ''love'', ''hear'', ''teamarchitect'', ''underlie'', ''package'', ''do'', ''work'', ''challenge'', ''visionocrhandwrite'', ''languagefinde'', ''contractsextracte'', ''clusteringworke'', ''require'', ''follow'', ''academia'', ''relate'', ''discriminate'', ''protect'', ''include'']'
$ws.Range("I7").Value = '[''go'', ''do'', ''use'', ''need'', ''patch'', ''look'', ''get'', ''enhance'', ''have'', ''complete'', ''realize'', ''propose'', ''include'', ''figure'', ''impact'', ''categorize'', ''join'', ''unite'', ''share'', ''make'', ''discover'', ''personalize'', ''provide'', ''consist'', ''develop'', ''purchase'', ''ecosystemfacilitate'', ''take'', ''check'', ''route'', ''send'', ''own'', ''bridge'', ''turn'', ''prioritize'', ''contribute'', ''put'', ''start'', ''create'', ''identify'', ''tour'', ''observe'', ''exist'', ''assist'', ''come'', ''includeworke'', ''talk'', ''schedule'', ''happen'', ''miss'', ''coordinate'', ''be'', ''imagine'', ''speak'', ''reach'', ''internshipbee'', ''orient'', ''organize'', ''pull'', ''involve'', ''enrol'', ''major'', ''demonstrate'', ''’s'', ''work'', ''establish'', ''explore'', ''implement'', ''invite'', ''participate'', ''network'', ''encourage'', ''bear'', ''discriminate'', ''protect'', ''prohibit'']'
